# edit.ps1 -- apply "change presentation to match proposal" edit
# Target: reorder/retitle the 9-slide "Attention Is All You Need" deck into
# the 16-slide "Brand Sentiment Analysis of Twitter Posts" deck.

$p = $ppt.ActivePresentation
$msoTrue = -1
$ctrAnchor = 3   # msoAnchorMiddle

# ---------------------------------------------------------------------
# STEP 1: Duplicate slide 1 (original title slide) -- the duplicate will
# become the NEW title slide ("Brand Sentiment Analysis of Twitter
# Posts" + the three author textboxes), landing right after slide 1.
# ---------------------------------------------------------------------
$dupRange = $p.Slides.Item(1).Duplicate()
$newTitleSlide = $dupRange.Item(1)

# -- Resize / retitle the duplicate's title placeholder --
$dTitle = $newTitleSlide.Shapes.Item(1)
$dTitle.Left = 72
$dTitle.Top = 162
$dTitle.Width = 816
$dTitle.Height = 216
$dTitle.TextFrame.TextRange.Text = "Brand Sentiment Analysis of Twitter Posts"
$dTitle.TextFrame.TextRange.Font.Size = 36

# -- Reposition the three author textboxes a bit lower & center them --
$dAuthor1 = $newTitleSlide.Shapes.Item(2)
$dAuthor1.Top = 324
$dAuthor1.TextFrame.VerticalAnchor = $ctrAnchor

$dAuthor2 = $newTitleSlide.Shapes.Item(3)
$dAuthor2.Left = 390.3958
$dAuthor2.Top = 324
$dAuthor2.TextFrame.VerticalAnchor = $ctrAnchor

$dAuthor3 = $newTitleSlide.Shapes.Item(4)
$dAuthor3.Left = 708.962
$dAuthor3.Top = 324
$dAuthor3.TextFrame.VerticalAnchor = $ctrAnchor

# ---------------------------------------------------------------------
# STEP 2: Turn the ORIGINAL slide 1 into the new intro-videoclip slide:
# drop the three author textboxes, retitle the remaining title shape.
# ---------------------------------------------------------------------
$introSlide = $p.Slides.Item(1)
# delete the 3 textboxes (Jonathan Agustin / Fernando Calderon / Juliet Lawton)
$introSlide.Shapes.Item(4).Delete()
$introSlide.Shapes.Item(3).Delete()
$introSlide.Shapes.Item(2).Delete()
$introSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Intro Videoclip & Segue to Presentation"

# ---------------------------------------------------------------------
# STEP 3: Delete the old "Related Work" slide (no longer used).
# ---------------------------------------------------------------------
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    if ($s.Shapes.Item(1).TextFrame.TextRange.Text -eq "Related Work") {
        $s.Delete()
        break
    }
}

# ---------------------------------------------------------------------
# STEP 4: Retitle the remaining original divider slides to match the
# new proposal outline (content/order handled later via MoveTo).
# ---------------------------------------------------------------------
function Set-TitleByOldText($oldText, $newText) {
    for ($i = 1; $i -le $p.Slides.Count; $i++) {
        $s = $p.Slides.Item($i)
        if ($s.Shapes.Item(1).TextFrame.TextRange.Text -eq $oldText) {
            $s.Shapes.Item(1).TextFrame.TextRange.Text = $newText
            return $s
        }
    }
    return $null
}

Set-TitleByOldText "Introduction" "Introduction" | Out-Null
$slideProblem  = Set-TitleByOldText "Problem Definition & AI Techniques" "Problem Statement"
$slideScope    = Set-TitleByOldText "Dataset Description" "Scope"
$slideAlgos    = Set-TitleByOldText "Experimental Design" "Algorithms"
$slideDatasets = Set-TitleByOldText "Results & Discussion" "Datasets"
$slideExpDsgn  = Set-TitleByOldText "Conclusion & Future Work" "Experimental Design"
$slideAssess   = Set-TitleByOldText "Acknowledgements" "First Stage: Assess"

# ---------------------------------------------------------------------
# STEP 5: Add the new stage / closing slides (all title-only, same
# "2 - Title and Content" layout used by the other divider slides).
# ---------------------------------------------------------------------
$layout2 = $p.SlideMaster.CustomLayouts.Item(2)

function Add-TitleOnlySlide($title) {
    $idx = $p.Slides.Count + 1
    $s = $p.Slides.AddSlide($idx, $layout2)
    # drop the unused body/content placeholder that comes from the layout
    for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
        $shp = $s.Shapes.Item($i)
        if ($shp.PlaceholderFormat.Type -ne 13) {
            if ($shp.Name -notlike "Title*") {
                $shp.Delete()
            }
        }
    }
    $s.Shapes.Item(1).TextFrame.TextRange.Text = $title
    return $s
}

$slideDevelop1 = Add-TitleOnlySlide "Second Stage: Develop"
$slideDevelop2 = Add-TitleOnlySlide "Second Stage: Develop"
$slideEvaluate = Add-TitleOnlySlide "Third Stage: Evaluate"
$slideResults  = Add-TitleOnlySlide "Results"
$slideChalleng = Add-TitleOnlySlide "Challenges"
$slideConclusion = Add-TitleOnlySlide "Conclusion & Future Work"
$slideAcknowledgements = Add-TitleOnlySlide "Acknowledgements"

# ---------------------------------------------------------------------
# STEP 6: Reorder everything into the final sequence.
# ---------------------------------------------------------------------
function MoveTitleTo($title, $pos) {
    for ($i = 1; $i -le $p.Slides.Count; $i++) {
        $s = $p.Slides.Item($i)
        if ($s.Shapes.Item(1).TextFrame.TextRange.Text -eq $title) {
            $s.MoveTo($pos)
            return
        }
    }
}

# Desired final order (1-based):
#  1 Intro Videoclip & Segue to Presentation
#  2 Brand Sentiment Analysis of Twitter Posts
#  3 Introduction
#  4 Problem Statement
#  5 Scope
#  6 Algorithms
#  7 Datasets
#  8 Experimental Design
#  9 First Stage: Assess
# 10 Second Stage: Develop
# 11 Second Stage: Develop
# 12 Third Stage: Evaluate
# 13 Results
# 14 Challenges
# 15 Conclusion & Future Work
# 16 Acknowledgements

$introSlide.MoveTo(1)
$newTitleSlide.MoveTo(2)
MoveTitleTo "Introduction" 3
$slideProblem.MoveTo(4)
$slideScope.MoveTo(5)
$slideAlgos.MoveTo(6)
$slideDatasets.MoveTo(7)
$slideExpDsgn.MoveTo(8)
$slideAssess.MoveTo(9)
$slideDevelop1.MoveTo(10)
$slideDevelop2.MoveTo(11)
$slideEvaluate.MoveTo(12)
$slideResults.MoveTo(13)
$slideChalleng.MoveTo(14)
$slideConclusion.MoveTo(15)
$slideAcknowledgements.MoveTo(16)

# ---------------------------------------------------------------------
# STEP 7: Set the speaker notes for every slide (final content).
# ---------------------------------------------------------------------
function Set-Notes($slide, $text) {
    $np = $slide.NotesPage
    for ($i = 1; $i -le $np.Shapes.Count; $i++) {
        $shp = $np.Shapes.Item($i)
        if ($shp.PlaceholderFormat.Type -eq 2) {
            $shp.TextFrame.TextRange.Text = $text
            return
        }
    }
}

Set-Notes $p.Slides.Item(1) "A cool introductory videoclip will be inserted here to WOW the audience and show that we mean business."
Set-Notes $p.Slides.Item(2) ("Let" + [char]0x2019 + "s get started. Welcome to our project titled " + [char]0x201C + "Brand Sentiment Analysis of Twitter Posts" + [char]0x201D + ". ")
Set-Notes $p.Slides.Item(3) "Our team has worked diligently to develop a system that leverages AI to understand public sentiment towards brands based on social media content."
Set-Notes $p.Slides.Item(4) "We build two models: a Brand Classifier that predicts whether a Twitter post expresses a brand, and a Brand Sentiment Analyzer that predicts the sentiment of a Twitter post towards a brand."
Set-Notes $p.Slides.Item(5) "They are trained on Twitter posts, so they cannot classify posts from other social media sites. They also process only text data, ignoring non-textual elements like images, videos, and audio clips."
Set-Notes $p.Slides.Item(6) ("We experimented with various machine learning algorithms like Na" + [char]0x00EF + "ve Bayes, Logistic Regression, Support Vector Machines, Recurrent Neural Networks, and Transformers.")
Set-Notes $p.Slides.Item(7) "We used the Sentiment140 dataset, which contains 1.6 million tweets, to train our models."
Set-Notes $p.Slides.Item(8) "Our methods were designed to ensure a thorough and systematic approach to the problem. It involved using at least two different types of AI and machine learning algorithms. We conducted an investigation of the analytics solution to the problem, which included aspects of experimental comparison. We also explored variable importance to understand which features were most important in making good predictions."
Set-Notes $p.Slides.Item(9) "In the assess stage, we analyzed the dataset and selected appropriate models based on the specific requirements of the sentiment analysis task, the characteristics of the available data, and the nature of the problem."
Set-Notes $p.Slides.Item(10) "In the development stage, we handled invalid values and discarded non-textual content. We trained our models to identify a brand and align with the sentiment expressed in a tweet, and then tuned them to improve performance."
Set-Notes $p.Slides.Item(11) "In the development stage, we handled invalid values and discarded non-textual content. We trained our models to identify a brand and align with the sentiment expressed in a tweet, and then tuned them to improve performance."
Set-Notes $p.Slides.Item(12) "In the evaluation stage, we measured the models' ability to identify brands and predict sentiment in Twitter posts. We used standard classification metrics such as Accuracy, Precision, Recall, and F1 Score to evaluate the models."
Set-Notes $p.Slides.Item(13) ("We present the results of our experiments, including any relevant figures or tables, discuss the implications of our results, and explain whether our results support our initial hypothesis. We present the results of our models. We achieved an accuracy rate over 50%, demonstrating the potential of machine learning in understanding public sentiment towards brands based on social media content.")
Set-Notes $p.Slides.Item(14) ("We faced challenges such as the " + [char]0x2018 + "aboutness" + [char]0x2019 + " problem, which refers to the challenge of determining the subject of the sentiment expressed in a sentence. We also encountered issues with complex sentences and ambiguous subjects.")
Set-Notes $p.Slides.Item(15) "In the conclusion and future work section, we will summarize the main findings of our project, conclude the report by discussing the significance of our findings, and discuss potential future work, such as how our project could be extended or improved."
Set-Notes $p.Slides.Item(16) "We would like to take a moment to express our deepest gratitude to those who have made this project possible. First and foremost, we would like to thank our professor, whose guidance and expertise have been invaluable throughout this process. Your patience and dedication have not gone unnoticed, and we are truly grateful for your support."

Write-Host "Done. Slide count: $($p.Slides.Count)"
